# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rows 16-46 hold one "Periodo Mora" record per row (periods 1..31). The
# period labels in column E are re-sorted from descending (2003 .. 1705)
# to ascending (1705 .. 2003), and the "Salario Basico" (F) / "Valor Mora"
# (G) figures that go with each period are refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for E16:E46 (was descending before the edit).
$periods = @(
    "1705","1706","1709","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$firstRow = 16
$lastRow  = 46

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i

    # Column E: period label for this row.
    $ws.Cells.Item($row, 5).Value = $periods[$i]

    # Column F: Salario Basico. Periods through 1808 (rows 16-27) use
    # 29509; periods from 1809 onward (rows 28-46) use 31249.
    if ($row -le 27) {
        $ws.Cells.Item($row, 6).Value = 29509
    } else {
        $ws.Cells.Item($row, 6).Value = 31249
    }

    # Column G: Valor Mora, refreshed to 781242 for every data row.
    $ws.Cells.Item($row, 7).Value = 781242
}
